$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall
$ws.Cells.Item(28, 8).Value = 969.8333
$ws.Cells.Item(28, 9).Value = 992.82355
$ws.Cells.Item(28, 10).Value = 914
$ws.Cells.Item(28, 11).Value = 992.82355
$ws.Cells.Item(28, 12).Value = 914
$ws.Cells.Item(28, 13).Value = -507.82355
$ws.Cells.Item(28, 14).Value = -1884

# Row 41: The Write Stuff
$ws.Cells.Item(41, 8).Value = 334
$ws.Cells.Item(41, 9).Value = 450.2
$ws.Cells.Item(41, 10).Value = 251
$ws.Cells.Item(41, 11).Value = 450.2
$ws.Cells.Item(41, 12).Value = 251
$ws.Cells.Item(41, 13).Value = -10.19999999999999
$ws.Cells.Item(41, 14).Value = -1131

# Row 74: Adhesive of Antipathy
$ws.Cells.Item(74, 8).Value = 3900
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 3900
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 3900
$ws.Cells.Item(74, 13).Value = ""
$ws.Cells.Item(74, 14).Value = -5772

# Row 77: It's Gonna Grow Back (L)
$ws.Cells.Item(77, 8).Value = 3900
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 3900
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 19500
$ws.Cells.Item(77, 13).Value = ""
$ws.Cells.Item(77, 14).Value = -28860

# Row 116: Growing Up
$ws.Cells.Item(116, 8).Value = 1730.6666
$ws.Cells.Item(116, 9).Value = 1595.25
$ws.Cells.Item(116, 10).Value = 2001.5
$ws.Cells.Item(116, 11).Value = 1595.25
$ws.Cells.Item(116, 12).Value = 2001.5
$ws.Cells.Item(116, 13).Value = 1846.75
$ws.Cells.Item(116, 14).Value = -8885.5

# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 2900895.8
$ws.Cells.Item(138, 10).Value = 5408487
$ws.Cells.Item(138, 12).Value = 16225461
$ws.Cells.Item(138, 14).Value = -16235741

# Row 141: Remedy for Reason
$ws.Cells.Item(141, 8).Value = 1863.6451
$ws.Cells.Item(141, 9).Value = 1571.3103
$ws.Cells.Item(141, 10).Value = 6102.5
$ws.Cells.Item(141, 11).Value = 4713.9309
$ws.Cells.Item(141, 12).Value = 18307.5
$ws.Cells.Item(141, 13).Value = 466.0690999999997
$ws.Cells.Item(141, 14).Value = -28667.5

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth
$ws.Cells.Item(5, 8).Value = 323.16666
$ws.Cells.Item(5, 9).Value = 229.8
$ws.Cells.Item(5, 10).Value = 790
$ws.Cells.Item(5, 11).Value = 229.8
$ws.Cells.Item(5, 12).Value = 790
$ws.Cells.Item(5, 13).Value = -117.8
$ws.Cells.Item(5, 14).Value = -1014

# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 4201.61
$ws.Cells.Item(32, 9).Value = 3766.957
$ws.Cells.Item(32, 10).Value = 9976.286
$ws.Cells.Item(32, 11).Value = 3766.957
$ws.Cells.Item(32, 12).Value = 9976.286
$ws.Cells.Item(32, 13).Value = -3479.957
$ws.Cells.Item(32, 14).Value = -10550.286

# Row 45: Hollow Hallmarks
$ws.Cells.Item(45, 8).Value = 1284.25
$ws.Cells.Item(45, 9).Value = 1251.6154
$ws.Cells.Item(45, 11).Value = 1251.6154
$ws.Cells.Item(45, 13).Value = -874.6153999999999

# Row 61: Dealing with the Tough Stuff
$ws.Cells.Item(61, 8).Value = 50101404
$ws.Cells.Item(61, 9).Value = 55612616
$ws.Cells.Item(61, 10).Value = 500500
$ws.Cells.Item(61, 11).Value = 55612616
$ws.Cells.Item(61, 12).Value = 500500
$ws.Cells.Item(61, 13).Value = -55612404
$ws.Cells.Item(61, 14).Value = -500924

# Row 69: The Cut Alembical Cord
$ws.Cells.Item(69, 8).Value = 53166.668
$ws.Cells.Item(69, 10).Value = 53166.668
$ws.Cells.Item(69, 12).Value = 53166.668
$ws.Cells.Item(69, 14).Value = -54664.668

# Row 72: Sheer Distill Power (L)
$ws.Cells.Item(72, 8).Value = 53166.668
$ws.Cells.Item(72, 10).Value = 53166.668
$ws.Cells.Item(72, 12).Value = 159500.004
$ws.Cells.Item(72, 14).Value = -166988.004

# Row 97: Ore for Me
$ws.Cells.Item(97, 8).Value = 2977089.5
$ws.Cells.Item(97, 9).Value = 3290098.8
$ws.Cells.Item(97, 10).Value = 3500
$ws.Cells.Item(97, 11).Value = 3290098.8
$ws.Cells.Item(97, 12).Value = 3500
$ws.Cells.Item(97, 13).Value = -3289602.8
$ws.Cells.Item(97, 14).Value = -4492

# Row 110: Scheduled Maintenance
$ws.Cells.Item(110, 8).Value = 2501831.5
$ws.Cells.Item(110, 9).Value = 10000000
$ws.Cells.Item(110, 10).Value = 2442
$ws.Cells.Item(110, 11).Value = 10000000
$ws.Cells.Item(110, 12).Value = 2442
$ws.Cells.Item(110, 13).Value = -9997955
$ws.Cells.Item(110, 14).Value = -6532

# Row 122: Haste for High Durium
$ws.Cells.Item(122, 8).Value = 4117375.5
$ws.Cells.Item(122, 9).Value = 2006.2778
$ws.Cells.Item(122, 10).Value = 12348114
$ws.Cells.Item(122, 11).Value = 6018.8334
$ws.Cells.Item(122, 12).Value = 37044342
$ws.Cells.Item(122, 13).Value = -3568.8334
$ws.Cells.Item(122, 14).Value = -37049242

# Row 136: Metal with Mettle
$ws.Cells.Item(136, 8).Value = 50101404
$ws.Cells.Item(136, 9).Value = 55612616
$ws.Cells.Item(136, 10).Value = 500500
$ws.Cells.Item(136, 11).Value = 166837848
$ws.Cells.Item(136, 12).Value = 1501500
$ws.Cells.Item(136, 13).Value = -166835298
$ws.Cells.Item(136, 14).Value = -1506600

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences
$ws.Cells.Item(4, 8).Value = 323.16666
$ws.Cells.Item(4, 9).Value = 229.8
$ws.Cells.Item(4, 10).Value = 790
$ws.Cells.Item(4, 11).Value = 229.8
$ws.Cells.Item(4, 12).Value = 790
$ws.Cells.Item(4, 13).Value = -114.8
$ws.Cells.Item(4, 14).Value = -1020

# Row 22: Riveting Run
$ws.Cells.Item(22, 8).Value = 234.83333
$ws.Cells.Item(22, 9).Value = 241.8
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(22, 11).Value = 241.8
$ws.Cells.Item(22, 12).Value = 200
$ws.Cells.Item(22, 13).Value = -68.80000000000001
$ws.Cells.Item(22, 14).Value = -546

# Row 94: High Steal
$ws.Cells.Item(94, 8).Value = 461
$ws.Cells.Item(94, 9).Value = 359.7
$ws.Cells.Item(94, 10).Value = 663.6
$ws.Cells.Item(94, 11).Value = 359.7
$ws.Cells.Item(94, 12).Value = 663.6
$ws.Cells.Item(94, 13).Value = 91.30000000000001
$ws.Cells.Item(94, 14).Value = -1565.6

# Row 99: Meddle in Metal
$ws.Cells.Item(99, 8).Value = 1636.25
$ws.Cells.Item(99, 9).Value = 1881.6666
$ws.Cells.Item(99, 10).Value = 900
$ws.Cells.Item(99, 11).Value = 1881.6666
$ws.Cells.Item(99, 12).Value = 900
$ws.Cells.Item(99, 13).Value = -383.6666
$ws.Cells.Item(99, 14).Value = -3896

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent
$ws.Cells.Item(7, 8).Value = 255.07692
$ws.Cells.Item(7, 9).Value = 121.5
$ws.Cells.Item(7, 10).Value = 314.44446
$ws.Cells.Item(7, 11).Value = 121.5
$ws.Cells.Item(7, 12).Value = 314.44446
$ws.Cells.Item(7, 13).Value = -8.5
$ws.Cells.Item(7, 14).Value = -540.4444599999999

# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 46249.914
$ws.Cells.Item(132, 9).Value = 2912.3635
$ws.Cells.Item(132, 10).Value = 85976
$ws.Cells.Item(132, 11).Value = 8737.0905
$ws.Cells.Item(132, 12).Value = 257928
$ws.Cells.Item(132, 13).Value = -6207.0905
$ws.Cells.Item(132, 14).Value = -262988

# Row 134: Wood You Be Quiet
$ws.Cells.Item(134, 8).Value = 144073.14
$ws.Cells.Item(134, 9).Value = 1418.6666
$ws.Cells.Item(134, 10).Value = 1000000
$ws.Cells.Item(134, 11).Value = 4255.9998
$ws.Cells.Item(134, 12).Value = 3000000
$ws.Cells.Item(134, 13).Value = -1720.9998
$ws.Cells.Item(134, 14).Value = -3005070

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Cells.Item(5, 8).Value = 835.61536
$ws.Cells.Item(5, 9).Value = 585.6111
$ws.Cells.Item(5, 10).Value = 1398.125
$ws.Cells.Item(5, 11).Value = 1756.8333
$ws.Cells.Item(5, 12).Value = 4194.375
$ws.Cells.Item(5, 13).Value = -1644.8333
$ws.Cells.Item(5, 14).Value = -4418.375

# Row 135: Not-so-secret Ingredient
$ws.Cells.Item(135, 8).Value = 835.61536
$ws.Cells.Item(135, 9).Value = 585.6111
$ws.Cells.Item(135, 10).Value = 1398.125
$ws.Cells.Item(135, 11).Value = 5270.4999
$ws.Cells.Item(135, 12).Value = 12583.125
$ws.Cells.Item(135, 13).Value = -2735.4999
$ws.Cells.Item(135, 14).Value = -17653.125

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Cells.Item(80, 8).Value = 3100.9443
$ws.Cells.Item(80, 9).Value = 2211.3333
$ws.Cells.Item(80, 10).Value = 3990.5557
$ws.Cells.Item(80, 11).Value = 2211.3333
$ws.Cells.Item(80, 12).Value = 3990.5557
$ws.Cells.Item(80, 13).Value = -1213.3333
$ws.Cells.Item(80, 14).Value = -5986.5557

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Cells.Item(83, 8).Value = 3100.9443
$ws.Cells.Item(83, 9).Value = 2211.3333
$ws.Cells.Item(83, 10).Value = 3990.5557
$ws.Cells.Item(83, 11).Value = 11056.6665
$ws.Cells.Item(83, 12).Value = 19952.7785
$ws.Cells.Item(83, 13).Value = -6064.666499999999
$ws.Cells.Item(83, 14).Value = -29936.7785

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Cells.Item(97, 8).Value = 1853.2
$ws.Cells.Item(97, 9).Value = 1790
$ws.Cells.Item(97, 10).Value = 2000.6666
$ws.Cells.Item(97, 11).Value = 1790
$ws.Cells.Item(97, 12).Value = 2000.6666
$ws.Cells.Item(97, 13).Value = -1294
$ws.Cells.Item(97, 14).Value = -2992.6666

# Row 102: Put the Metal to the Peddle
$ws.Cells.Item(102, 8).Value = 2365.6191
$ws.Cells.Item(102, 9).Value = 2347.75
$ws.Cells.Item(102, 11).Value = 2347.75
$ws.Cells.Item(102, 13).Value = -725.75

# Row 107: Whetstones for the Workers
$ws.Cells.Item(107, 8).Value = 380.41666
$ws.Cells.Item(107, 9).Value = 223.04347
$ws.Cells.Item(107, 10).Value = 4000
$ws.Cells.Item(107, 11).Value = 223.04347
$ws.Cells.Item(107, 12).Value = 4000
$ws.Cells.Item(107, 13).Value = 1696.95653
$ws.Cells.Item(107, 14).Value = -7840

# Row 122: Awarding Academic Excellence
$ws.Cells.Item(122, 8).Value = 1994.8667
$ws.Cells.Item(122, 9).Value = 1717.6364
$ws.Cells.Item(122, 10).Value = 2757.25
$ws.Cells.Item(122, 11).Value = 5152.9092
$ws.Cells.Item(122, 12).Value = 8271.75
$ws.Cells.Item(122, 13).Value = -2702.9092
$ws.Cells.Item(122, 14).Value = -13171.75

# Row 135: Fan of the Foreign
$ws.Cells.Item(135, 8).Value = 26968.234
$ws.Cells.Item(135, 10).Value = 26968.234
$ws.Cells.Item(135, 12).Value = 26968.234
$ws.Cells.Item(135, 14).Value = -37108.234

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Cells.Item(40, 8).Value = 2701.25
$ws.Cells.Item(40, 9).Value = 1800
$ws.Cells.Item(40, 10).Value = 3001.6667
$ws.Cells.Item(40, 11).Value = 1800
$ws.Cells.Item(40, 12).Value = 3001.6667
$ws.Cells.Item(40, 13).Value = -1664
$ws.Cells.Item(40, 14).Value = -3273.6667

# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 31788.652
$ws.Cells.Item(136, 9).Value = 18562.174
$ws.Cells.Item(136, 10).Value = 170666.67
$ws.Cells.Item(136, 11).Value = 55686.522
$ws.Cells.Item(136, 12).Value = 512000.01
$ws.Cells.Item(136, 13).Value = -53136.522
$ws.Cells.Item(136, 14).Value = -517100.01

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke
$ws.Cells.Item(62, 8).Value = 4002.5
$ws.Cells.Item(62, 9).Value = 4002
$ws.Cells.Item(62, 10).Value = 4003
$ws.Cells.Item(62, 11).Value = 4002
$ws.Cells.Item(62, 12).Value = 4003
$ws.Cells.Item(62, 13).Value = -3378
$ws.Cells.Item(62, 14).Value = -5251

# Row 65: Desperate for Diversionaries (L)
$ws.Cells.Item(65, 8).Value = 4002.5
$ws.Cells.Item(65, 9).Value = 4002
$ws.Cells.Item(65, 10).Value = 4003
$ws.Cells.Item(65, 11).Value = 20010
$ws.Cells.Item(65, 12).Value = 20015
$ws.Cells.Item(65, 13).Value = -16890
$ws.Cells.Item(65, 14).Value = -26255

# Row 96: Skills on Display
$ws.Cells.Item(96, 8).Value = 9996
$ws.Cells.Item(96, 9).Value = 3740
$ws.Cells.Item(96, 10).Value = 14166.667
$ws.Cells.Item(96, 11).Value = 3740
$ws.Cells.Item(96, 12).Value = 14166.667
$ws.Cells.Item(96, 13).Value = -2367
$ws.Cells.Item(96, 14).Value = -16912.667

# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 40803.39
$ws.Cells.Item(132, 9).Value = 35533.62
$ws.Cells.Item(132, 10).Value = 47749.91
$ws.Cells.Item(132, 11).Value = 106600.86
$ws.Cells.Item(132, 12).Value = 143249.73
$ws.Cells.Item(132, 13).Value = -104070.86
$ws.Cells.Item(132, 14).Value = -148309.73
